$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 3 and row 4 (date, volume, min/max/avg price columns)
$ws.Range("D3").Value = 44993
$ws.Range("M3").Value = 14
$ws.Range("O3").Value = 200000
$ws.Range("P3").Value = 190000
$ws.Range("S3").Value = 190000

$ws.Range("D4").Value = 44672
$ws.Range("M4").Value = 8
$ws.Range("O4").Value = 180000
$ws.Range("P4").Value = 180000
$ws.Range("S4").Value = 180000
